# The published diff shows that every data row from the old row 49 down to the
# old row 146 has been shifted down by exactly one row (old row 49 -> new row 50,
# ... old row 146 -> new row 147), and a brand new weekly record has been
# inserted at row 49. The sheet's dimension grows from A1:R146 to A1:R147.
#
# Inserting a whole row above the old row 49 reproduces this shift (and the
# accompanying dimension/format bookkeeping) automatically, so we only need to
# insert the row and then fill in the values for the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 49:146 down to 50:147, leaving a blank row 49 in their place.
$ws.Rows("49").Insert()

# Populate the new row 49 with the new "Feria Lagunitas de Puerto Montt -
# Arveja Verde" record added by this edit.
$ws.Range("A49").Value2 = 4
$ws.Range("B49").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C49").Value2 = "Los Lagos"
$ws.Range("D49").Value2 = 44883
$ws.Range("E49").Value2 = 10
$ws.Range("F49").Value2 = 100112022
$ws.Range("G49").Value2 = "Arveja Verde"
$ws.Range("H49").Value2 = "Sin especificar"
$ws.Range("I49").Value2 = "Primera"
$ws.Range("J49").Value2 = 80
$ws.Range("K49").Value2 = 27000
$ws.Range("L49").Value2 = 27000
$ws.Range("M49").Value2 = 27000
$ws.Range("N49").Value2 = "`$/saco 25 kilos"
$ws.Range("O49").Value2 = "Región del Maule"
$ws.Range("P49").Value2 = 1080
$ws.Range("Q49").Value2 = 25
$ws.Range("R49").Value2 = "Hortaliza"
